# Fix data summary report:
#  - give the two mixed-species tabs underscore-safe names
#  - re-autofit column C (WHON5_CODE) on the sheets where it was left too
#    narrow, which also moves the "active sheet" bookmark around as a side
#    effect of selecting those columns
#  - end with the Acinetobacter tab active/selected, matching the final
#    on-screen state the user left the workbook in

$wb = $excel.ActiveWorkbook

$acineto = $wb.Worksheets.Item("Acinetobacter species")
$pseudo  = $wb.Worksheets.Item("Pseudomonas aeruginosa")

$acineto.Name = "Acinetobacter_species"
$pseudo.Name  = "Pseudomonas_aeruginosa"

# Pseudomonas_aeruginosa: column C (WHON5_CODE) was too narrow - select it
# and autofit it to its content, same as column C already is on the
# Acinetobacter_species sheet.
$pseudo.Columns("C").Select()
$pseudo.Columns("C").EntireColumn.AutoFit()
$pseudo.Columns("C").ColumnWidth = 13.1666666666667

# Same fix on ENTEROBACTERIACEAE_X_SAL_SHI.
$entero = $wb.Worksheets.Item("ENTEROBACTERIACEAE_X_SAL_SHI")
$entero.Columns("C").Select()
$entero.Columns("C").EntireColumn.AutoFit()
$entero.Columns("C").ColumnWidth = 13.1666666666667

# Leave the workbook with the Acinetobacter tab active, cell L9 selected.
$acineto.Activate()
$acineto.Range("L9").Select()
